$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "26.161.27"
$rng.Style = "Normal"
$rng = $ws.Range("E2")
$rng.NumberFormat = "@"
$rng.Value = "  -1.19%  "
$rng.Style = "Normal"
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "1.658.31"
$rng.Style = "Normal"
$rng = $ws.Range("E3")
$rng.NumberFormat = "@"
$rng.Value = "  -1.20%  "
$rng.Style = "Normal"
$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$rng.Value = "  +0.36%  "
$rng.Style = "Normal"
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "216.09"
$rng.Style = "Normal"
$rng = $ws.Range("E5")
$rng.NumberFormat = "@"
$rng.Value = "  -1.63%  "
$rng.Style = "Normal"
$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = "0.5206"
$rng.Style = "Normal"
$rng = $ws.Range("E6")
$rng.NumberFormat = "@"
$rng.Value = "  -2.60%  "
$rng.Style = "Normal"
$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = "0.2627"
$rng.Style = "Normal"
$rng = $ws.Range("E8")
$rng.NumberFormat = "@"
$rng.Value = "  -3.24%  "
$rng.Style = "Normal"
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = "0.06261"
$rng.Style = "Normal"
$rng = $ws.Range("E9")
$rng.NumberFormat = "@"
$rng.Value = "  -2.47%  "
$rng.Style = "Normal"
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "20.73"
$rng.Style = "Normal"
$rng = $ws.Range("E10")
$rng.NumberFormat = "@"
$rng.Value = "  -5.57%  "
$rng.Style = "Normal"
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "0.07710"
$rng.Style = "Normal"
$rng = $ws.Range("E11")
$rng.NumberFormat = "@"
$rng.Value = "  -1.15%  "
$rng.Style = "Normal"
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = "1.653.64"
$rng.Style = "Normal"
$rng = $ws.Range("E12")
$rng.NumberFormat = "@"
$rng.Value = "  -1.51%  "
$rng.Style = "Normal"
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "4.418"
$rng.Style = "Normal"
$rng = $ws.Range("E13")
$rng.NumberFormat = "@"
$rng.Value = "  -2.24%  "
$rng.Style = "Normal"
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "1.883.58"
$rng.Style = "Normal"
$rng = $ws.Range("E14")
$rng.NumberFormat = "@"
$rng.Value = "  -1.30%  "
$rng.Style = "Normal"
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "0.5421"
$rng.Style = "Normal"
$rng = $ws.Range("E15")
$rng.NumberFormat = "@"
$rng.Value = "  -3.23%  "
$rng.Style = "Normal"
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = "0.0₅8124"
$rng.Style = "Normal"
$rng = $ws.Range("E16")
$rng.NumberFormat = "@"
$rng.Value = "  -2.68%  "
$rng.Style = "Normal"
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "64.38"
$rng.Style = "Normal"
$rng = $ws.Range("E17")
$rng.NumberFormat = "@"
$rng.Value = "  -2.21%  "
$rng.Style = "Normal"
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "26.189.92"
$rng.Style = "Normal"
$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$rng.Value = "  -1.27%  "
$rng.Style = "Normal"
$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$rng.Value = "  +0.44%  "
$rng.Style = "Normal"
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "4.624"
$rng.Style = "Normal"
$rng = $ws.Range("E20")
$rng.NumberFormat = "@"
$rng.Value = "  -3.89%  "
$rng.Style = "Normal"
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "191.68"
$rng.Style = "Normal"
$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$rng.Value = "  -0.99%  "
$rng.Style = "Normal"
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "10.05"
$rng.Style = "Normal"
$rng = $ws.Range("E22")
$rng.NumberFormat = "@"
$rng.Value = "  -2.73%  "
$rng.Style = "Normal"
$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = "6.057"
$rng.Style = "Normal"
$rng = $ws.Range("E23")
$rng.NumberFormat = "@"
$rng.Value = "  -4.29%  "
$rng.Style = "Normal"
$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$rng.Value = "  +0.46%  "
$rng.Style = "Normal"
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "139.61"
$rng.Style = "Normal"
$rng = $ws.Range("E25")
$rng.NumberFormat = "@"
$rng.Value = "  -1.92%  "
$rng.Style = "Normal"
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "0.1225"
$rng.Style = "Normal"
$rng = $ws.Range("E26")
$rng.NumberFormat = "@"
$rng.Value = "  -4.53%  "
$rng.Style = "Normal"
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "7.173"
$rng.Style = "Normal"
$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$rng.Value = "  -3.36%  "
$rng.Style = "Normal"
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "16.04"
$rng.Style = "Normal"
$rng = $ws.Range("E28")
$rng.NumberFormat = "@"
$rng.Value = "  -1.78%  "
$rng.Style = "Normal"
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "1.412"
$rng.Style = "Normal"
$rng = $ws.Range("E29")
$rng.NumberFormat = "@"
$rng.Value = "  -2.21%  "
$rng.Style = "Normal"
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = "0.05988"
$rng.Style = "Normal"
$rng = $ws.Range("E30")
$rng.NumberFormat = "@"
$rng.Value = "  -5.10%  "
$rng.Style = "Normal"
$rng = $ws.Range("E31")
$rng.NumberFormat = "@"
$rng.Value = "  -1.26%  "
$rng.Style = "Normal"
$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = "3.547"
$rng.Style = "Normal"
$rng = $ws.Range("E32")
$rng.NumberFormat = "@"
$rng.Value = "  -1.76%  "
$rng.Style = "Normal"
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "3.235"
$rng.Style = "Normal"
$rng = $ws.Range("E33")
$rng.NumberFormat = "@"
$rng.Value = "  -6.63%  "
$rng.Style = "Normal"
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "1.611"
$rng.Style = "Normal"
$rng = $ws.Range("E34")
$rng.NumberFormat = "@"
$rng.Value = "  -5.29%  "
$rng.Style = "Normal"
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = "0.9657"
$rng.Style = "Normal"
$rng = $ws.Range("E35")
$rng.NumberFormat = "@"
$rng.Value = "  -4.69%  "
$rng.Style = "Normal"
$rng = $ws.Range("D36")
$rng.NumberFormat = "@"
$rng.Value = "2.417"
$rng.Style = "Normal"
$rng = $ws.Range("E36")
$rng.NumberFormat = "@"
$rng.Value = "  -0.11%  "
$rng.Style = "Normal"
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = "2.769"
$rng.Style = "Normal"
$rng = $ws.Range("E37")
$rng.NumberFormat = "@"
$rng.Value = "  -0.61%  "
$rng.Style = "Normal"
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = "0.5667"
$rng.Style = "Normal"
$rng = $ws.Range("E38")
$rng.NumberFormat = "@"
$rng.Value = "  -7.87%  "
$rng.Style = "Normal"
$rng = $ws.Range("B39")
$rng.NumberFormat = "@"
$rng.Value = "VeChain"
$rng.Style = "Normal"
$rng = $ws.Range("C39")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$rng.Style = "Normal"
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = "0.01598"
$rng.Style = "Normal"
$rng = $ws.Range("E39")
$rng.NumberFormat = "@"
$rng.Value = "  -2.30%  "
$rng.Style = "Normal"
$rng = $ws.Range("B40")
$rng.NumberFormat = "@"
$rng.Value = "FraxShare"
$rng.Style = "Normal"
$rng = $ws.Range("C40")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$rng.Style = "Normal"
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "6.002"
$rng.Style = "Normal"
$rng = $ws.Range("E40")
$rng.NumberFormat = "@"
$rng.Value = "  -2.50%  "
$rng.Style = "Normal"
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = "0.8549"
$rng.Style = "Normal"
$rng = $ws.Range("E41")
$rng.NumberFormat = "@"
$rng.Value = "  -1.40%  "
$rng.Style = "Normal"
$rng = $ws.Range("E42")
$rng.NumberFormat = "@"
$rng.Value = "  +0.39%  "
$rng.Style = "Normal"
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = "1.011.94"
$rng.Style = "Normal"
$rng = $ws.Range("E43")
$rng.NumberFormat = "@"
$rng.Value = "  -7.39%  "
$rng.Style = "Normal"
$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = "100.17"
$rng.Style = "Normal"
$rng = $ws.Range("E44")
$rng.NumberFormat = "@"
$rng.Value = "  -0.47%  "
$rng.Style = "Normal"
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = "1.799.22"
$rng.Style = "Normal"
$rng = $ws.Range("E45")
$rng.NumberFormat = "@"
$rng.Value = "  -1.34%  "
$rng.Style = "Normal"
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "0.0₈111"
$rng.Style = "Normal"
$rng = $ws.Range("E46")
$rng.NumberFormat = "@"
$rng.Value = "  -1.48%  "
$rng.Style = "Normal"
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "56.74"
$rng.Style = "Normal"
$rng = $ws.Range("E47")
$rng.NumberFormat = "@"
$rng.Value = "  -3.65%  "
$rng.Style = "Normal"
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "1.007"
$rng.Style = "Normal"
$rng = $ws.Range("E48")
$rng.NumberFormat = "@"
$rng.Value = "  +0.54%  "
$rng.Style = "Normal"
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = "7.972"
$rng.Style = "Normal"
$rng = $ws.Range("E49")
$rng.NumberFormat = "@"
$rng.Value = "  -2.71%  "
$rng.Style = "Normal"
$rng = $ws.Range("E50")
$rng.NumberFormat = "@"
$rng.Value = "  -0.74%  "
$rng.Style = "Normal"
$rng = $ws.Range("B51")
$rng.NumberFormat = "@"
$rng.Value = "RenderToken"
$rng.Style = "Normal"
$rng = $ws.Range("C51")
$rng.NumberFormat = "@"
$rng.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$rng.Style = "Normal"
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "1.452"
$rng.Style = "Normal"
$rng = $ws.Range("E51")
$rng.NumberFormat = "@"
$rng.Value = "  -1.54%  "
$rng.Style = "Normal"
